$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "52×27=1404" "56×32=1792"
Replace-Text "17×72=1224" "68×90=6120"
Replace-Text "93×68=6324" "27×79=2133"
Replace-Text "91×50=4550" "85×95=8075"
Replace-Text "26×66=1716" "31×43=1333"
Replace-Text "56×95=5320" "68×40=2720"
Replace-Text "69×26=1794" "13×94=1222"
Replace-Text "25×62=1550" "99×87=8613"
Replace-Text "94×92=8648" "43×22=946"
Replace-Text "17×17=289" "71×16=1136"
Replace-Text "73×91=6643" "20×99=1980"
Replace-Text "21×16=336" "35×87=3045"
Replace-Text "13×21=273" "43×56=2408"
Replace-Text "18×70=1260" "58×95=5510"
Replace-Text "54×18=972" "16×42=672"
Replace-Text "59×56=3304" "30×63=1890"
Replace-Text "90×93=8370" "68×14=952"
Replace-Text "12×82=984" "33×61=2013"
Replace-Text "69×75=5175" "29×13=377"
Replace-Text "65×26=1690" "80×23=1840"
Replace-Text "76×38=2888" "89×30=2670"
Replace-Text "52×93=4836" "90×93=8370"
Replace-Text "76×82=6232" "30×59=1770"
Replace-Text "85×87=7395" "14×83=1162"
Replace-Text "90×62=5580" "16×74=1184"
